# "updated main GSC export data"
#
# The "Chart" sheet holds a daily time-series table (Date | Non-HTTPS URLs |
# HTTPS URLs) in A1:C91. This edit rolls the window forward by one day:
# the oldest date (2025-10-08, row 2) is dropped, every remaining row's
# Date/HTTPS-URLs values shifts up by one row, and a new row for
# 2026-01-06 is appended at the end (row 91) with its HTTPS URL count.
# The "Non-HTTPS URLs" column (B) stays 0 throughout and is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Date column text for A2:A91, in row order.
$dates = @("2025-10-09","2025-10-10","2025-10-11","2025-10-12","2025-10-13","2025-10-14","2025-10-15","2025-10-16","2025-10-17","2025-10-18","2025-10-19","2025-10-20","2025-10-21","2025-10-22","2025-10-23","2025-10-24","2025-10-25","2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18","2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27","2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06","2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11","2025-12-12","2025-12-13","2025-12-14","2025-12-15","2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23","2025-12-24","2025-12-25","2025-12-26","2025-12-27","2025-12-28","2025-12-29","2025-12-30","2025-12-31","2026-01-01","2026-01-02","2026-01-03","2026-01-04","2026-01-05","2026-01-06")

# New "HTTPS URLs" counts for C2:C91, in row order (aligned with $dates).
$counts = @(11.0,16.0,22.0,33.0,43.0,51.0,58.0,58.0,68.0,66.0,71.0,71.0,75.0,82.0,92.0,84.0,85.0,86.0,90.0,83.0,90.0,93.0,92.0,102.0,113.0,115.0,107.0,105.0,100.0,94.0,86.0,83.0,66.0,54.0,46.0,43.0,40.0,37.0,35.0,30.0,29.0,26.0,25.0,25.0,26.0,26.0,25.0,25.0,27.0,28.0,28.0,27.0,27.0,27.0,27.0,27.0,26.0,25.0,25.0,25.0,26.0,27.0,27.0,29.0,29.0,30.0,30.0,31.0,31.0,31.0,31.0,31.0,32.0,32.0,32.0,32.0,30.0,31.0,32.0,30.0,28.0,28.0,28.0,28.0,29.0,29.0,28.0,27.0,27.0,28.0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 2 + $i

    # Column A holds plain text dates (shared strings), not real Excel
    # dates. Assigning a date-look-alike string straight to .Value makes
    # Excel "smart" parse it into a date serial number and stamp the cell
    # with a date NumberFormat/style. Pre-formatting the cell as Text before
    # the assignment keeps it a literal string; ClearFormats() afterwards
    # drops that temporary Text format again so the cell's style stays the
    # default (style 0), exactly as it was originally.
    $aCell = $ws.Range("A" + $r)
    $aCell.NumberFormat = "@"
    $aCell.Value = $dates[$i]
    $aCell.ClearFormats()

    # Column C is a plain number, no special handling needed.
    $ws.Range("C" + $r).Value = $counts[$i]
}
